# Append a new movie row (Goodfellas / Martin Scorsese / 148) to Sheet1,
# mirroring a user typing a new record at the end of the table and then
# leaving the selection on the next empty row (A5) — same as Excel does
# after Enter-ing through a row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "Goodfellas"
$ws.Range("B4").Value = "Martin Scorsese"
$ws.Range("C4").Value = 148

$ws.Range("A5").Select()
